$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 33 (ALC, G=5512)
$ws_ALC.Range("H33").Value = 823.3200000000001
$ws_ALC.Range("I33").Value = 877.5217
$ws_ALC.Range("K33").Value = 877.5217
$ws_ALC.Range("M33").Value = -648.5217

# Row 40 (ALC, G=5505)
$ws_ALC.Range("H40").Value = 1907.3
$ws_ALC.Range("I40").Value = 1957.6923
$ws_ALC.Range("J40").Value = 1813.7142
$ws_ALC.Range("K40").Value = 1957.6923
$ws_ALC.Range("L40").Value = 1813.7142
$ws_ALC.Range("M40").Value = -1782.6923
$ws_ALC.Range("N40").Value = -2163.7142

# Row 43 (ALC, G=5472)
$ws_ALC.Range("H43").Value = 644.9583
$ws_ALC.Range("I43").Value = 510.5
$ws_ALC.Range("J43").Value = 712.1875
$ws_ALC.Range("K43").Value = 510.5
$ws_ALC.Range("L43").Value = 712.1875
$ws_ALC.Range("M43").Value = -441.5
$ws_ALC.Range("N43").Value = -850.1875

# Row 132 (ALC, G=44049)
$ws_ALC.Range("H132").Value = 4067.7556
$ws_ALC.Range("I132").Value = 4497.9644
$ws_ALC.Range("J132").Value = 3359.1765
$ws_ALC.Range("K132").Value = 13493.8932
$ws_ALC.Range("L132").Value = 10077.5295
$ws_ALC.Range("M132").Value = -10963.8932
$ws_ALC.Range("N132").Value = -15137.5295

# Row 107 (BSM, G=27706)
$ws_BSM.Range("H107").Value = 1259.375
$ws_BSM.Range("I107").Value = 966.3889
$ws_BSM.Range("J107").Value = 2138.3333
$ws_BSM.Range("K107").Value = 966.3889
$ws_BSM.Range("L107").Value = 2138.3333
$ws_BSM.Range("M107").Value = 953.6111
$ws_BSM.Range("N107").Value = -5978.3333

# Row 134 (BSM, G=43998)
$ws_BSM.Range("H134").Value = 1646.3793
$ws_BSM.Range("I134").Value = 1260.2084
$ws_BSM.Range("J134").Value = 3500
$ws_BSM.Range("K134").Value = 3780.6252
$ws_BSM.Range("L134").Value = 10500
$ws_BSM.Range("M134").Value = -1245.6252
$ws_BSM.Range("N134").Value = -15570

# Row 58 (CRP, G=44021)
$ws_CRP.Range("H58").Value = 2520.5483
$ws_CRP.Range("I58").Value = 899.3333
$ws_CRP.Range("J58").Value = 4765.3076
$ws_CRP.Range("K58").Value = 899.3333
$ws_CRP.Range("L58").Value = 4765.3076
$ws_CRP.Range("M58").Value = -696.3333
$ws_CRP.Range("N58").Value = -5171.3076

# Row 99 (CRP, G=36198)
$ws_CRP.Range("H99").Value = 2105267.5
$ws_CRP.Range("I99").Value = 3251484.8
$ws_CRP.Range("J99").Value = 3869
$ws_CRP.Range("K99").Value = 3251484.8
$ws_CRP.Range("L99").Value = 3869
$ws_CRP.Range("M99").Value = -3249986.8
$ws_CRP.Range("N99").Value = -6865

# Row 126 (CRP, G=36198)
$ws_CRP.Range("H126").Value = 2105267.5
$ws_CRP.Range("I126").Value = 3251484.8
$ws_CRP.Range("J126").Value = 3869
$ws_CRP.Range("K126").Value = 9754454.399999999
$ws_CRP.Range("L126").Value = 11607
$ws_CRP.Range("M126").Value = -9751984.399999999
$ws_CRP.Range("N126").Value = -16547

# Row 134 (CRP, G=44020)
$ws_CRP.Range("H134").Value = 3708.4614
$ws_CRP.Range("I134").Value = 4930.467
$ws_CRP.Range("J134").Value = 2042.091
$ws_CRP.Range("K134").Value = 14791.401
$ws_CRP.Range("L134").Value = 6126.272999999999
$ws_CRP.Range("M134").Value = -12256.401
$ws_CRP.Range("N134").Value = -11196.273

# Row 136 (CRP, G=44021)
$ws_CRP.Range("H136").Value = 2520.5483
$ws_CRP.Range("I136").Value = 899.3333
$ws_CRP.Range("J136").Value = 4765.3076
$ws_CRP.Range("K136").Value = 2697.9999
$ws_CRP.Range("L136").Value = 14295.9228
$ws_CRP.Range("M136").Value = -147.9998999999998
$ws_CRP.Range("N136").Value = -19395.9228

# Row 2 (CUL, G=4847)
$ws_CUL.Range("H2").Value = 165094.25
$ws_CUL.Range("I2").Value = 330051.84
$ws_CUL.Range("J2").Value = 136.66667
$ws_CUL.Range("K2").Value = 1980311.04
$ws_CUL.Range("L2").Value = 820.0000200000001
$ws_CUL.Range("M2").Value = -1980198.04
$ws_CUL.Range("N2").Value = -1046.00002

# Row 5 (CUL, G=43974)
$ws_CUL.Range("H5").Value = 1141.0769
$ws_CUL.Range("J5").Value = 1525
$ws_CUL.Range("L5").Value = 4575
$ws_CUL.Range("N5").Value = -4799

# Row 38 (CUL, G=4860)
$ws_CUL.Range("H38").Value = 244
$ws_CUL.Range("I38").Value = 241.8
$ws_CUL.Range("J38").Value = 246.2
$ws_CUL.Range("K38").Value = 725.4000000000001
$ws_CUL.Range("L38").Value = 738.5999999999999
$ws_CUL.Range("M38").Value = -378.4000000000001
$ws_CUL.Range("N38").Value = -1432.6

# Row 135 (CUL, G=43974)
$ws_CUL.Range("H135").Value = 1141.0769
$ws_CUL.Range("J135").Value = 1525
$ws_CUL.Range("L135").Value = 13725
$ws_CUL.Range("N135").Value = -18795

# Row 43 (GSM, G=4218)
$ws_GSM.Range("H43").Value = 1047.091
$ws_GSM.Range("I43").Value = 1047.091
$ws_GSM.Range("K43").Value = 1047.091
$ws_GSM.Range("M43").Value = -896.0909999999999

# Row 132 (GSM, G=44008)
$ws_GSM.Range("H132").Value = 2039.2222
$ws_GSM.Range("I132").Value = 1223.5
$ws_GSM.Range("J132").Value = 3670.6667
$ws_GSM.Range("K132").Value = 3670.5
$ws_GSM.Range("L132").Value = 11012.0001
$ws_GSM.Range("M132").Value = -1140.5
$ws_GSM.Range("N132").Value = -16072.0001

# Row 16 (LTW, G=5289)
$ws_LTW.Range("H16").Value = 665.8823
$ws_LTW.Range("I16").Value = 612.6667
$ws_LTW.Range("J16").Value = 1065
$ws_LTW.Range("K16").Value = 612.6667
$ws_LTW.Range("L16").Value = 1065
$ws_LTW.Range("M16").Value = -442.6667
$ws_LTW.Range("N16").Value = -1405

# Row 64 (LTW, G=10810)
$ws_LTW.Range("H64").Value = 23630
$ws_LTW.Range("I64").Value = 8000
$ws_LTW.Range("J64").Value = 27537.5
$ws_LTW.Range("K64").Value = 8000
$ws_LTW.Range("L64").Value = 27537.5
$ws_LTW.Range("M64").Value = -7775
$ws_LTW.Range("N64").Value = -27987.5

# Row 67 (LTW, G=10810)
$ws_LTW.Range("H67").Value = 23630
$ws_LTW.Range("I67").Value = 8000
$ws_LTW.Range("J67").Value = 27537.5
$ws_LTW.Range("K67").Value = 8000
$ws_LTW.Range("L67").Value = 27537.5
$ws_LTW.Range("M67").Value = -7220
$ws_LTW.Range("N67").Value = -29097.5

# Row 68 (LTW, G=12563)
$ws_LTW.Range("H68").Value = 11668595
$ws_LTW.Range("I68").Value = 56390124
$ws_LTW.Range("J68").Value = 2108.4783
$ws_LTW.Range("K68").Value = 56390124
$ws_LTW.Range("L68").Value = 2108.4783
$ws_LTW.Range("M68").Value = -56389375
$ws_LTW.Range("N68").Value = -3606.4783

# Row 71 (LTW, G=12563)
$ws_LTW.Range("H71").Value = 11668595
$ws_LTW.Range("I71").Value = 56390124
$ws_LTW.Range("J71").Value = 2108.4783
$ws_LTW.Range("K71").Value = 281950620
$ws_LTW.Range("L71").Value = 10542.3915
$ws_LTW.Range("M71").Value = -281946876
$ws_LTW.Range("N71").Value = -18030.3915

# Row 122 (LTW, G=36247)
$ws_LTW.Range("H122").Value = 3852.4
$ws_LTW.Range("I122").Value = 4376
$ws_LTW.Range("J122").Value = 3503.3333
$ws_LTW.Range("K122").Value = 13128
$ws_LTW.Range("L122").Value = 10509.9999
$ws_LTW.Range("M122").Value = -10678
$ws_LTW.Range("N122").Value = -15409.9999

# Row 132 (LTW, G=44058)
$ws_LTW.Range("H132").Value = 10006382
$ws_LTW.Range("I132").Value = 27791116
$ws_LTW.Range("J132").Value = 2468.4062
$ws_LTW.Range("K132").Value = 83373348
$ws_LTW.Range("L132").Value = 7405.2186
$ws_LTW.Range("M132").Value = -83370818
$ws_LTW.Range("N132").Value = -12465.2186

# Row 15 (WVR, G=2670)
$ws_WVR.Range("H15").Value = 7658.6665
$ws_WVR.Range("J15").Value = 7658.6665
$ws_WVR.Range("L15").Value = 7658.6665
$ws_WVR.Range("N15").Value = -8234.666499999999

# Row 62 (WVR, G=12589)
$ws_WVR.Range("H62").Value = 5357.143
$ws_WVR.Range("I62").Value = 5540
$ws_WVR.Range("J62").Value = 4900
$ws_WVR.Range("K62").Value = 5540
$ws_WVR.Range("L62").Value = 4900
$ws_WVR.Range("M62").Value = -4916
$ws_WVR.Range("N62").Value = -6148

# Row 63 (WVR, G=10824)
$ws_WVR.Range("H63").Value = 30249
$ws_WVR.Range("J63").Value = 30249
$ws_WVR.Range("L63").Value = 30249
$ws_WVR.Range("N63").Value = -31497

# Row 65 (WVR, G=12589)
$ws_WVR.Range("H65").Value = 5357.143
$ws_WVR.Range("I65").Value = 5540
$ws_WVR.Range("J65").Value = 4900
$ws_WVR.Range("K65").Value = 27700
$ws_WVR.Range("L65").Value = 24500
$ws_WVR.Range("M65").Value = -24580
$ws_WVR.Range("N65").Value = -30740

# Row 66 (WVR, G=10824)
$ws_WVR.Range("H66").Value = 30249
$ws_WVR.Range("J66").Value = 30249
$ws_WVR.Range("L66").Value = 90747
$ws_WVR.Range("N66").Value = -96987

# Row 122 (WVR, G=36208)
$ws_WVR.Range("H122").Value = 1000
$ws_WVR.Range("I122").Value = 0
$ws_WVR.Range("J122").Value = 1000
$ws_WVR.Range("K122").Value = 0
$ws_WVR.Range("L122").Value = 3000
$ws_WVR.Range("N122").Value = -7900
$ws_WVR.Range("M122").ClearContents()

# Row 126 (WVR, G=36210)
$ws_WVR.Range("H126").Value = 3290.1875
$ws_WVR.Range("I126").Value = 3588.7856
$ws_WVR.Range("J126").Value = 1200
$ws_WVR.Range("K126").Value = 10766.3568
$ws_WVR.Range("L126").Value = 3600
$ws_WVR.Range("M126").Value = -8296.356800000001
$ws_WVR.Range("N126").Value = -8540

# Row 132 (WVR, G=44029)
$ws_WVR.Range("H132").Value = 1445.75
$ws_WVR.Range("I132").Value = 754.96
$ws_WVR.Range("J132").Value = 2597.0667
$ws_WVR.Range("K132").Value = 2264.88
$ws_WVR.Range("L132").Value = 7791.2001
$ws_WVR.Range("M132").Value = 265.1199999999999
$ws_WVR.Range("N132").Value = -12851.2001

# Row 136 (WVR, G=44031)
$ws_WVR.Range("H136").Value = 1769.8889
$ws_WVR.Range("I136").Value = 1094
$ws_WVR.Range("J136").Value = 2200
$ws_WVR.Range("K136").Value = 3282
$ws_WVR.Range("L136").Value = 6600
$ws_WVR.Range("N136").Value = -11700
$ws_WVR.Range("M136").Value = -732
